$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B54 to be a numeric value instead of text
$ws.Range("B54").Value = 4

# Add new row 55 with data
$ws.Range("A55").Value = "Sunsi Wu"
$ws.Range("B55").Value = "3"
$ws.Range("C55").Value = "无"
$ws.Range("D55").Value = "ACK"
$ws.Range("E55").Value = "EXP"
$ws.Range("F55").Value = "0ffe4b07-d72b-4753-8576-ca80ee89bdb3"
$ws.Range("G55").Value = "SJzMATlAZ_annotated.xlsx"
$ws.Range("H55").Value = "We avoid using k-means because it requires knowing the number of clusters a priory."
